# "Add files via upload" — the workbook data was re-uploaded with an
# updated value in column B (Track.xlsx, sheet "data").
# Functional change: InstructorID for TrackID 4 (row 5) goes from 4 to 8,
# and the active selection in the sheet moves to that cell (B5).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$ws.Range("B5").Value = 8
$ws.Range("B5").Select() | Out-Null
